# Update eggs and weights
# Appends the two new daily tracking rows (2/11/2021 and 2/12/2021) to the
# bottom of the "chinups"/"pushups" log on Sheet1, and moves the viewport /
# selection the way the author left the workbook (scrolled down, with H31
# selected) before saving.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Append new data rows -------------------------------------------------
# Row 23 (the prior last row) holds date serial 44237 (2/10/2021) in column A
# with chinups (B) / pushups (C) counts. Continue the sequence for two more
# days.
$ws.Cells.Item(24, 1).Value = 44238   # 2/11/2021
$ws.Cells.Item(24, 2).Value = 0
$ws.Cells.Item(24, 3).Value = 0

$ws.Cells.Item(25, 1).Value = 44239   # 2/12/2021
$ws.Cells.Item(25, 2).Value = 50
$ws.Cells.Item(25, 3).Value = 45

# Copy the date formatting (m/d/yyyy) from the existing column A cells onto
# the two new date cells so the new rows match the rest of the column.
$ws.Range("A23").Copy()
$ws.Range("A24:A25").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

# --- Update the view state -------------------------------------------------
# Scroll the window so row 16 is at the top and select H31, matching where
# the author left the sheet.
$win = $excel.ActiveWindow
$win.ScrollRow = 16
$win.ScrollColumn = 1
$ws.Range("H31").Select()
